$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15 ("The Bahea") status moved from "Edited" to "On IG"
$ws.Range("D15").Value = "On IG"

# New photo rows appended to the schedule
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Plane In the Clouds"
$ws.Range("D16").Value = "Shot"

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Sunset"
$ws.Range("D17").Value = "Shot"

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Ausie Goose"
$ws.Range("D18").Value = "Shot"

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Turtle on A Log"
$ws.Range("D19").Value = "Shot"

$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "Squrille"
$ws.Range("D20").Value = "Shot"

# Match the new selection left behind by the edit
$ws.Range("D16").Select()
